# Fruta / hortaliza, semanal
# A new weekly price record is inserted as row 26 (pushing the existing
# rows 26-36 down to 27-37), for "Agrícola del Norte S.A. de Arica" -
# Poroto verde, Magnum variety, sourced from Perú.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26; this shifts existing rows 26-36 -> 27-37
# and preserves formatting/styles from the surrounding rows (e.g. the date
# format style on column D).
$ws.Rows.Item(26).Insert()

$ws.Cells.Item(26, 1).Value = 1
$ws.Cells.Item(26, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(26, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(26, 4).Value = 44435
$ws.Cells.Item(26, 5).Value = 15
$ws.Cells.Item(26, 6).Value = 100112031
$ws.Cells.Item(26, 7).Value = 'Poroto verde'
$ws.Cells.Item(26, 8).Value = 'Magnum'
$ws.Cells.Item(26, 9).Value = 'Primera'
$ws.Cells.Item(26, 10).Value = 160
$ws.Cells.Item(26, 11).Value = 21000
$ws.Cells.Item(26, 12).Value = 22000
$ws.Cells.Item(26, 13).Value = 21500
$ws.Cells.Item(26, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(26, 15).Value = 'Perú'
$ws.Cells.Item(26, 16).Value = 860
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = 'Hortaliza'

# Make sure the date cell keeps the same numeric date format used by the
# rest of the "Fecha" column.
$ws.Cells.Item(26, 4).NumberFormat = $ws.Cells.Item(27, 4).NumberFormat
